# update to manual status column
# - Rows 3 & 4 (samples 2 & 3): manualStatus column I switches from the
#   numeric placeholder "4" to the textual status "[4]"
# - Column F (fastqFileName) is widened so the long file names are readable
# - Active selection moves to I5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "[4]"
$ws.Range("I4").Value = "[4]"

# Target stored width is 56.96 characters; Excel quantizes column widths to
# whole pixels (MDW=6px), so 56.16667 is the input that rounds to the same
# pixel width (342px -> 57.0) closest to the recorded 56.96.
$ws.Columns.Item(6).ColumnWidth = 56.16666666666667

$ws.Range("I5").Select()
